# Adds the "Materiales y suministros" / "Servicios" detail rows (11-15) to
# the "Gasto Capital" sheet, per the report-preview commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gasto Capital")
$xlPasteValues = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues

# Helper: write a numeric-looking value as real TEXT (shared string), not a
# number, without leaving any residual formula or cell style behind. We bake
# a TEXT() formula into a literal value via copy / paste-special-values.
function Set-TextValue {
    param($range, [string]$text)
    $range.Formula = '=TEXT("' + $text + '","@")'
    $range.Copy() | Out-Null
    $range.PasteSpecial($xlPasteValues) | Out-Null
    $excel.CutCopyMode = $false
}

# --- Row 11: "Materiales y suministros" section header -------------------
$ws.Range("B11").Value = "Materiales y suministros"

# --- Row 12: FILLER DE 0.0040 A 0.009 25 HOJAS X 4 ------------------------
$ws.Range("B12").Value = "FILLER DE 0.0040 A 0.009 25 HOJAS X 4"
$ws.Range("C12").Value = 1
Set-TextValue $ws.Range("D12") "400716"
$ws.Range("E12").Value = 58.98

# --- Row 13: FILLER DE 0.050 A 1.000  20 HOJAS X 12  ----------------------
$ws.Range("B13").Value = "FILLER DE 0.050 A 1.000  20 HOJAS X 12 "
$ws.Range("C13").Value = 1
Set-TextValue $ws.Range("D13") "400716"
$ws.Range("E13").Value = 81.42
$ws.Range("F13").Value = "X"
$ws.Range("G13").Value = "X"
$ws.Range("N13").Value = "X"

# --- Row 14: "Servicios" section header -----------------------------------
$ws.Range("B14").Value = "Servicios"

# --- Row 15: MANTENIMIENTO DE EXTINTORES ----------------------------------
$ws.Range("B15").Value = "MANTENIMIENTO DE EXTINTORES"
$ws.Range("C15").Value = 0
$ws.Range("E15").Value = 0
